$wb = $excel.ActiveWorkbook

# --- Sheet "Totales Junio": update selection to B7 ---
$wsJunio = $wb.Worksheets.Item("Totales Junio")
$wsJunio.Activate()
$wsJunio.Range("B7").Select()

# --- Sheet "Julio": fill in figures for July 7 and July 8, 2022 and update selection ---
$wsJulio = $wb.Worksheets.Item("Julio")
$wsJulio.Activate()

# Row 9 -> 7 de julio de 2022 (serial 44749)
$wsJulio.Range("B9").Value = 91
$wsJulio.Range("C9").Value = 3
$wsJulio.Range("D9").Value = 1
$wsJulio.Range("E9").Value = 0
$wsJulio.Range("F9").Value = 0
$wsJulio.Range("G9").Value = 0
$wsJulio.Range("H9").Value = 0
$wsJulio.Range("I9").Value = 0
$wsJulio.Range("J9").Value = 1

# Row 10 -> 8 de julio de 2022 (serial 44750)
$wsJulio.Range("B10").Value = 114
$wsJulio.Range("C10").Value = 3
$wsJulio.Range("D10").Value = 4
$wsJulio.Range("E10").Value = 0
$wsJulio.Range("F10").Value = 1
$wsJulio.Range("G10").Value = 0
$wsJulio.Range("H10").Value = 0
$wsJulio.Range("I10").Value = 0
$wsJulio.Range("J10").Value = 0

$wsJulio.Range("K11").Select()
